$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the project placeholder text (shared string used by A1)
$ws.Range("A1").Value = "EXAMPLE PROJECT"

# 2. Widen column A so the longer project name is readable
$ws.Columns("A").ColumnWidth = 19

# 3. A new 3D model (.STL part) was added to the project -- fill in its
#    row of measurements (piece count, time, material) as row 9, right
#    after the existing data (rows 10+ keep their positions).
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 112.50800323486328
$ws.Range("D9").Value = 3
